# Updates workbook "北京-漫展信息.xlsx" to the state generated at commit 456a3b4.
# Touches three sheets:
#   - 展览 (exhibitions): inserts one new record (row 33), shifting the
#     remaining 2024-12-14.. rows down, and bumps several "want to go" counts.
#   - 本地生活 (local life): bumps two "want to go" counts.
#   - 全部类型 (all types): bumps several "want to go" counts (independent
#     running counters from 展览 / 演出 / 本地生活, not re-sorted/re-inserted).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------------
$ws1 = $wb.Sheets.Item("展览")

# Simple "want to go" (F) count refreshes for existing rows (no shift).
$ws1.Range("F4").Value  = 5890
$ws1.Range("F7").Value  = 539
$ws1.Range("F9").Value  = 1569
$ws1.Range("F13").Value = 1583
$ws1.Range("F14").Value = 1583
$ws1.Range("F15").Value = 1550
$ws1.Range("F17").Value = 148
$ws1.Range("F19").Value = 4420
$ws1.Range("F20").Value = 30
$ws1.Range("F23").Value = 815
$ws1.Range("F24").Value = 7
$ws1.Range("F25").Value = 46
$ws1.Range("F26").Value = 2304

# Insert a brand-new listing before row 33 (everything from the old row 33
# down shifts to row+1; a new row 37 appears at the bottom of the table).
$ws1.Rows("33").Insert()

# The freshly inserted row 33 loses the shared "index column" style that the
# rest of column A carries (s="1", bold + border + centered) -- restore it by
# copying the format from the cell directly below, which still has it.
$ws1.Range("A34").Copy($ws1.Range("A33"))

# Column A is a plain positional index (row number - 1), not content that
# travels with a row -- native Insert() shifted the old numbers down along
# with everything else, so re-stamp A33:A37 back to the expected sequence.
$ws1.Range("A33").Value = 32
$ws1.Range("A34").Value = 33
$ws1.Range("A35").Value = 34
$ws1.Range("A36").Value = 35
$ws1.Range("A37").Value = 36

# Date column (B) holds plain text like "2024-12-14", not a real Excel date;
# force text entry via a temporary Text number format, then drop back to the
# sheet's default (unstyled) look so the cell matches its neighbours.
$ws1.Range("B33").NumberFormat = "@"
$ws1.Range("B33").Value = "2024-12-14"
$ws1.Range("B33").Style = "Normal"

$ws1.Range("C33").Value = "北京·thebONE×GOJO超次元动漫游戏嘉年华17th"
$ws1.Range("D33").Value = "B1层西区会员活动中心（阳坊涮肉对面） 北投购物公园"
$ws1.Range("E33").Value = "2024.12.14 09:30-12.15 17:00"
$ws1.Range("F33").Value = 0
$ws1.Range("G33").Value = 1
$ws1.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=93495"
$ws1.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202410/Fzz24Usj1728969298701.jpeg"

# The shift carried the old row-34/35 "want to go" counts along verbatim;
# bump the two that also changed in this refresh.
$ws1.Range("F35").Value = 1206
$ws1.Range("F36").Value = 1191

# ---------------------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------------------
$ws3 = $wb.Sheets.Item("本地生活")
$ws3.Range("F3").Value = 661
$ws3.Range("F5").Value = 268

# ---------------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------------
$ws4 = $wb.Sheets.Item("全部类型")
$ws4.Range("F6").Value  = 661
$ws4.Range("F8").Value  = 5890
$ws4.Range("F16").Value = 539
$ws4.Range("F19").Value = 1569
$ws4.Range("F23").Value = 1583
$ws4.Range("F25").Value = 1550
$ws4.Range("F27").Value = 148
$ws4.Range("F29").Value = 4420
$ws4.Range("F32").Value = 815
$ws4.Range("F33").Value = 46
$ws4.Range("F35").Value = 2304
$ws4.Range("F47").Value = 1206
$ws4.Range("F49").Value = 1191
